$wb = $excel.ActiveWorkbook

# --- Architecture sheet: expand the "No error handling strategies" comment ---
$wsArch = $wb.Worksheets.Item("Architecture")
$wsArch.Range("F7").Value = "No error handling strategies on lower architectural layers. Error handling only on requests layer."

# --- Code sheet: row for check "C09" (There are confusion in use of the signs.) ---
# Move the checkmark from the YES column (C) to the NO column (D),
# and clear the outdated comment about == vs === in column F.
$wsCode = $wb.Worksheets.Item("Code")
$wsCode.Range("C11").Value = ""
$wsCode.Range("D11").Value = "X"
$wsCode.Range("F11").Value = ""

# --- Update the active selections / active sheet to match the saved view state ---
$wsReq = $wb.Worksheets.Item("Requirements")
$wsReq.Range("F17").Select()

$wsArch.Range("G7").Select()

$wsCode.Range("D11").Select()
$wsCode.Activate()
